$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on Hoja1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("✅ 1000 Bs = 7.69 = 31169.47 pesos", "✅ 1000 Bs = 7.52 = 30440.44 pesos")
$text = $text.Replace("✅ 31169.47 pesos = 7.67 = 962.31 Bs", "✅ 30440.44 pesos = 7.46 = 923.95 Bs")
$cell.Value = $text

# --- Update rate cells on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 133.05
$wsTasas.Range("O10").Value = 4050.1
$wsTasas.Range("N12").Value = 4082
$wsTasas.Range("O12").Value = 123.9
